{"js": "// The commit adds two new reference bullets (each consisting of some lead-in\n// text followed by a hyperlink) to the numbered/bulleted reference list,\n// right after the \"Showing / Hiding password in editText...\" item and right\n// before the \"Image Sources\" heading paragraph.\n//\n// (The rest of the underlying XML diff is made up of the Word spell/grammar\n// checker re-splitting already-existing runs and stamping <w:proofErr/>\n// markers around technical tokens such as \"Doxygen\", \"editText\", \"apk\",\n// \"android:shape\", etc. Those markers are produced internally by Word's\n// proofing pass; they carry no visible text or formatting change and there\n// is no Word JavaScript API call that authors them, so there is nothing to\n// replicate there - the paragraph text/content is identical before and\n// after. Likewise the shifted hyperlink relationship ids (rId13 -> rId15,\n// etc.) are just a side effect of the two new relationships being minted;\n// the actual hyperlink targets are unchanged.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the anchor paragraph (\"Showing / Hiding password in editText of\n// type textPassword <link>\") by its distinctive lead-in text instead of a\n// hard-coded index, so the script is resilient to minor paragraph-count\n// differences.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Showing / Hiding password in\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Could not find the 'Showing / Hiding password' paragraph\");\n}\n\n// --- New bullet 1: spinner text size/colour reference -------------------\nconst spinnerUrl =\n  \"https://stackoverflow.com/questions/9476665/how-to-change-spinner-text-size-and-text-color\";\nconst spinnerPara = anchor.insertParagraph(\n  \"Changing spinner text size and colour \",\n  Word.InsertLocation.after\n);\nconst spinnerLinkRange = spinnerPara.insertText(spinnerUrl, Word.InsertLocation.end);\nspinnerLinkRange.hyperlink = spinnerUrl;\n\n// --- New bullet 2: ratingBar colour reference ----------------------------\nconst ratingBarUrl =\n  \"https://stackoverflow.com/questions/32810341/android-change-color-of-ratingbar-to-golden\";\nconst ratingBarPara = spinnerPara.insertParagraph(\n  \"Changing the colour of the ratingBar \",\n  Word.InsertLocation.after\n);\nconst ratingBarLinkRange = ratingBarPara.insertText(ratingBarUrl, Word.InsertLocation.end);\nratingBarLinkRange.hyperlink = ratingBarUrl;\n\nawait context.sync();\n", "ps1": "# The commit adds two new reference bullets (each consisting of some lead-in\n# text followed by a hyperlink) to the numbered/bulleted reference list,\n# right after the \"Showing / Hiding password in editText...\" item and right\n# before the \"Image Sources\" heading paragraph.\n#\n# (The rest of the underlying XML diff is made up of the Word spell/grammar\n# checker re-splitting already-existing runs and stamping <w:proofErr/>\n# markers around technical tokens such as \"Doxygen\", \"editText\", \"apk\",\n# \"android:shape\", etc. Those markers are produced internally by Word's\n# proofing pass; they carry no visible text or formatting change and there\n# is no COM call that authors them, so there is nothing to replicate there -\n# the paragraph text/content is identical before and after. Likewise the\n# shifted hyperlink relationship ids (rId13 -> rId15, etc.) are just a side\n# effect of the two new relationships being minted; the actual hyperlink\n# targets are unchanged.)\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"Showing / Hiding password in editText of\n# type textPassword <link>\") by its distinctive lead-in text instead of a\n# hard-coded index, so the script is resilient to minor paragraph-count\n# differences.\n$anchor = $null\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"Showing / Hiding password in*\") {\n        $anchor = $p\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Could not find the 'Showing / Hiding password' paragraph\"\n}\n\n# --- New bullet 1: spinner text size/colour reference --------------------\n$spinnerLead = \"Changing spinner text size and colour \"\n$spinnerUrl = \"https://stackoverflow.com/questions/9476665/how-to-change-spinner-text-size-and-text-color\"\n\n$insertAt = $anchor.Range\n$insertAt.Collapse(0)          # wdCollapseEnd\n$insertAt.InsertParagraphAfter()\n\n$paras = $d.Paragraphs\n$spinnerPara = $paras.Item($anchor.Index + 1)\n$spinnerRange = $spinnerPara.Range\n$spinnerRange.Text = $spinnerLead + $spinnerUrl\n$fullRange = $spinnerPara.Range\n$urlStart = $fullRange.Start + $spinnerLead.Length\n$urlEnd = $fullRange.End - 1    # exclude the paragraph mark\n$urlRange = $d.Range($urlStart, $urlEnd)\n$d.Hyperlinks.Add($urlRange, $spinnerUrl) | Out-Null\n\n# --- New bullet 2: ratingBar colour reference -----------------------------\n$ratingLead = \"Changing the colour of the ratingBar \"\n$ratingUrl = \"https://stackoverflow.com/questions/32810341/android-change-color-of-ratingbar-to-golden\"\n\n$paras = $d.Paragraphs\n$spinnerPara = $paras.Item($anchor.Index + 1)\n$insertAt2 = $spinnerPara.Range\n$insertAt2.Collapse(0)         # wdCollapseEnd\n$insertAt2.InsertParagraphAfter()\n\n$paras = $d.Paragraphs\n$ratingPara = $paras.Item($anchor.Index + 2)\n$ratingRange = $ratingPara.Range\n$ratingRange.Text = $ratingLead + $ratingUrl\n$fullRange2 = $ratingPara.Range\n$urlStart2 = $fullRange2.Start + $ratingLead.Length\n$urlEnd2 = $fullRange2.End - 1  # exclude the paragraph mark\n$urlRange2 = $d.Range($urlStart2, $urlEnd2)\n$d.Hyperlinks.Add($urlRange2, $ratingUrl) | Out-Null\n"}
